$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1990.909
$ws.Range("I40").Value = 1590
$ws.Range("K40").Value = 1590
$ws.Range("M40").Value = -1415
$ws.Range("H125").Value = 337376.66
$ws.Range("I125").Value = 232
$ws.Range("K125").Value = 2088
$ws.Range("M125").Value = 372
$ws.Range("H127").Value = 2096.4285
$ws.Range("I127").Value = 2096.4285
$ws.Range("K127").Value = 6289.2855
$ws.Range("M127").Value = -1329.2855
$ws.Range("H137").Value = 4259.8945
$ws.Range("J137").Value = 4583.222
$ws.Range("L137").Value = 13749.666
$ws.Range("N137").Value = -18849.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4999
$ws.Range("J8").Value = 4999
$ws.Range("L8").Value = 4999
$ws.Range("N8").Value = -5287
$ws.Range("H12").Value = 7400
$ws.Range("I12").Value = 4001.5
$ws.Range("K12").Value = 4001.5
$ws.Range("M12").Value = -3828.5
$ws.Range("H16").Value = 659.3333
$ws.Range("I16").Value = 489
$ws.Range("K16").Value = 489
$ws.Range("M16").Value = -202
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H88").Value = 3778
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3778
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3778
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4590
$ws.Range("H91").Value = 3778
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3778
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3778
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6586
$ws.Range("H110").Value = 2779
$ws.Range("I110").Value = 2965
$ws.Range("K110").Value = 2965
$ws.Range("M110").Value = -920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17477.818
$ws.Range("H85").Value = 17477.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 663.0625
$ws.Range("I7").Value = 640.6
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 640.6
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -527.6
$ws.Range("N7").Value = -1226
$ws.Range("H58").Value = 2726
$ws.Range("I58").Value = 2650.3333
$ws.Range("K58").Value = 2650.3333
$ws.Range("M58").Value = -2447.3333
$ws.Range("H105").Value = 650
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 650
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 650
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -4144
$ws.Range("H132").Value = 2997.5
$ws.Range("I132").Value = 495
$ws.Range("K132").Value = 1485
$ws.Range("M132").Value = 1045
$ws.Range("H134").Value = 6188.3335
$ws.Range("I134").Value = 5426
$ws.Range("K134").Value = 16278
$ws.Range("M134").Value = -13743
$ws.Range("H136").Value = 2726
$ws.Range("I136").Value = 2650.3333
$ws.Range("K136").Value = 7950.999899999999
$ws.Range("M136").Value = -5400.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3592.182
$ws.Range("I97").Value = 2500.5715
$ws.Range("K97").Value = 2500.5715
$ws.Range("M97").Value = -2004.5715
$ws.Range("H113").Value = 1705.7142
$ws.Range("I113").Value = 1185.25
$ws.Range("J113").Value = 2399.6667
$ws.Range("K113").Value = 1185.25
$ws.Range("L113").Value = 2399.6667
$ws.Range("M113").Value = 984.75
$ws.Range("N113").Value = -6739.6667
$ws.Range("H122").Value = 17298.572
$ws.Range("I122").Value = 22218.8
$ws.Range("K122").Value = 66656.39999999999
$ws.Range("M122").Value = -64206.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 30125.5
$ws.Range("I68").Value = 3500.6667
$ws.Range("K68").Value = 3500.6667
$ws.Range("M68").Value = -2751.6667
$ws.Range("H71").Value = 30125.5
$ws.Range("I71").Value = 3500.6667
$ws.Range("K71").Value = 17503.3335
$ws.Range("M71").Value = -13759.3335
$ws.Range("H76").Value = 29355.5
$ws.Range("J76").Value = 29355.5
$ws.Range("L76").Value = 29355.5
$ws.Range("N76").Value = -30031.5
$ws.Range("H79").Value = 29355.5
$ws.Range("J79").Value = 29355.5
$ws.Range("L79").Value = 29355.5
$ws.Range("N79").Value = -31695.5
$ws.Range("H82").Value = 1312.4
$ws.Range("I82").Value = 1265.5
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 1265.5
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -904.5
$ws.Range("N82").Value = -2222
$ws.Range("H85").Value = 1312.4
$ws.Range("I85").Value = 1265.5
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 1265.5
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -17.5
$ws.Range("N85").Value = -3996
$ws.Range("H122").Value = 1045
$ws.Range("I122").Value = 1045
$ws.Range("K122").Value = 3135
$ws.Range("M122").Value = -685
$ws.Range("H132").Value = 4999.5
$ws.Range("I132").Value = 4999.5
$ws.Range("K132").Value = 14998.5
$ws.Range("M132").Value = -12468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H81").Value = 999
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937
$ws.Range("H84").Value = 999
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686
$ws.Range("H107").Value = 1454.6
$ws.Range("I107").Value = 943.25
$ws.Range("K107").Value = 2829.75
$ws.Range("M107").Value = -909.75
$ws.Range("H122").Value = 1749.75
$ws.Range("I122").Value = 1252
$ws.Range("K122").Value = 3756
$ws.Range("M122").Value = -1306
